$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$timestamp = "2025-11-07 01:19:40"

# Remove all existing hyperlinks up front; they will be re-created below so
# that every F-column cell ends up with a correct, freshly-bound hyperlink
# (avoids any stale ref/target mismatch from the row shuffle).
$ws.Hyperlinks.Delete()

# ---- Row 2 (unchanged content, timestamp refreshed) ----
$ws.Range("A2").Value = $timestamp
$ws.Range("B2").Value = "【Next.js × TypeScript × Tailwind】コンポーネント制作パートナー募集!"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5428507"
$ws.Range("G2").Value = 528
$ws.Range("H2").Value = "🔥AI,Next.js"

# ---- Row 3 (new posting inserted ahead of the rest) ----
$ws.Range("A3").Value = $timestamp
$ws.Range("B3").Value = "専門データ分析:AIコスト最適化設計と厳格な機密保持を必須とするWebシステム開発(段階的継続発注)"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5428695"
$ws.Range("G3").Value = 403
$ws.Range("H3").Value = "🔥AI,Ai ◆開発,システム開発"

# ---- Row 4 (was row 3) ----
$ws.Range("A4").Value = $timestamp
$ws.Range("B4").Value = "大手製造業向け センサー画像解析・高画質化のR&Dを支援するAIエンジニア募集(画像生成/超解像)"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5427956"
$ws.Range("G4").Value = 310
$ws.Range("H4").Value = "🔥AI,Ai"

# ---- Row 5 (was row 4) ----
$ws.Range("A5").Value = $timestamp
$ws.Range("B5").Value = "【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5217096"
$ws.Range("G5").Value = 243
$ws.Range("H5").Value = "🔥API ◆ツール"

# ---- Row 6 (was row 5) ----
$ws.Range("A6").Value = $timestamp
$ws.Range("B6").Value = "【急募】GitHub管理のBootstrapサイト移行作業依頼"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5428337"
$ws.Range("G6").Value = 58
$ws.Range("H6").Value = "◇サイト"

# ---- Row 7 (was row 6) ----
$ws.Range("A7").Value = $timestamp
$ws.Range("B7").Value = "勤怠管理システムの改修依頼"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5428278"
$ws.Range("G7").Value = 53
$ws.Range("H7").Value = "◇管理"

# ---- Row 8 (was row 7) ----
$ws.Range("A8").Value = $timestamp
$ws.Range("B8").Value = "【急募】WEB会計アプリ機能修正!納期11/09希望"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5428124"
$ws.Range("G8").Value = 38
$ws.Range("H8").Value = "◇アプリ"

# ---- Row 9 (new) ----
$ws.Range("A9").Value = $timestamp
$ws.Range("B9").Value = "【若手歓迎×リモートOK】SRE/インフラエンジニア(Google Cloud/長期・金融系案件)"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5428756"
$ws.Range("G9").Value = 25

# ---- Row 10 (new) ----
$ws.Range("A10").Value = $timestamp
$ws.Range("B10").Value = "【リーダー募集×リモートOK】SRE/インフラエンジニア(Google Cloud/長期金融系案件)"
$ws.Range("C10").Value = "システム開発"
$ws.Range("D10").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E10").Value = "期限情報なし"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5428755"
$ws.Range("G10").Value = 25

# ---- Row 11 (new) ----
$ws.Range("A11").Value = $timestamp
$ws.Range("B11").Value = "【急募】Googleworkスペース・ハブスポットのサーバー設定依頼"
$ws.Range("C11").Value = "システム開発"
$ws.Range("D11").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E11").Value = "期限情報なし"
$ws.Range("F11").Value = "https://www.lancers.jp/work/detail/5428509"
$ws.Range("G11").Value = 10

# Re-create the F-column hyperlinks (row order = final row order) so the
# relationship ids line up 1:1 with rows 2..11, and re-apply the Hyperlink
# character style that Excel normally stamps on a linked cell.
# (Literal URLs are used here rather than reading `.Value` back, since the
# target is known up-front from the same literals written above.)
$urls = @(
    "https://www.lancers.jp/work/detail/5428507",
    "https://www.lancers.jp/work/detail/5428695",
    "https://www.lancers.jp/work/detail/5427956",
    "https://www.lancers.jp/work/detail/5217096",
    "https://www.lancers.jp/work/detail/5428337",
    "https://www.lancers.jp/work/detail/5428278",
    "https://www.lancers.jp/work/detail/5428124",
    "https://www.lancers.jp/work/detail/5428756",
    "https://www.lancers.jp/work/detail/5428755",
    "https://www.lancers.jp/work/detail/5428509"
)
for ($i = 0; $i -lt $urls.Length; $i++) {
    $r = $i + 2
    $cell = $ws.Range("F" + $r)
    $ws.Hyperlinks.Add($cell, $urls[$i])
    $cell.Style = $ws.Range("F2").Style
}

# Column width tweaks (Excel's ColumnWidth is in "characters"; it renders
# ~0.8333 wider in the stored OOXML width because of the default 5px cell
# padding at this font/size, so back that padding out to land on 30 / 19).
$ws.Columns.Item(4).ColumnWidth = 29.166666666666668
$ws.Columns.Item(8).ColumnWidth = 18.166666666666668
